# wg test correction #2
# - Correct the AD5/AD6 input (0.3 -> 0.2); the dependent formulas
#   (AB5/AB6, AH5/AH6, AN5/AN6) recalc automatically off of this.
# - Update the sheet view: zoom 85% -> 130%, and the active
#   selection moves from F13 to AB6 (scrolled right toward column S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction ---------------------------------------------------
$ws.Range("AD5").Value = 0.2
$ws.Range("AD6").Value = 0.2

# --- View state ----------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 19          # topLeftCell column "S"
$win.ScrollRow = 1              # topLeftCell row 1
$win.Zoom = 130

# Move the selection/active cell to AB6 (also drives scroll position)
$ws.Range("AB6").Select()
